# Generate Report for Handoff
# Update the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# and mark the handoff type ("Priority" column) as "ht" for the files that were
# just (re-)handed off: 460d96af..., 4cfd3238..., 528a0249..., baa21ccc...,
# d29e39ad..., f04ac914... (rows 7, 8, 9, 12, 13, 14 in each per-language sheet).

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 12, 13, 14)

# --- Overview sheet: "Latest HO Xliff Generate Date" (column G) ---
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Cells.Item($r, 7).Value = "2016-08-30 14:25:29"
}

# --- zh-cn sheet: Priority (column E) + Latest Handoff Datetime (column H) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-30 14:25:19"
}

# --- de-de sheet: Priority (column E) + Latest Handoff Datetime (column H) ---
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-30 14:25:29"
}
